# Apply updated price/profit figures to each class sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below were recomputed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 9091319
$ws.Range("I38").Value = 10000151
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 30000453
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -30000081
$ws.Range("N38").Value = -9744
$ws.Range("H98").Value = 1299.25
$ws.Range("I98").Value = 1166
$ws.Range("K98").Value = 1166
$ws.Range("M98").Value = 332
$ws.Range("H122").Value = 1299.25
$ws.Range("I122").Value = 1166
$ws.Range("K122").Value = 3498
$ws.Range("M122").Value = -1048
$ws.Range("H137").Value = 1999
$ws.Range("I137").Value = 1998
$ws.Range("K137").Value = 5994
$ws.Range("M137").Value = -3444
$ws.Range("H138").Value = 4639.2593
$ws.Range("I138").Value = 2095.5557
$ws.Range("K138").Value = 6286.6671
$ws.Range("M138").Value = -1146.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3503363.5
$ws.Range("I32").Value = 3687926.5
$ws.Range("K32").Value = 3687926.5
$ws.Range("M32").Value = -3687639.5
$ws.Range("H61").Value = 1383.8334
$ws.Range("I61").Value = 1383.8334
$ws.Range("K61").Value = 1383.8334
$ws.Range("M61").Value = -1171.8334
$ws.Range("H74").Value = 4708
$ws.Range("I74").Value = 3849.6
$ws.Range("K74").Value = 3849.6
$ws.Range("M74").Value = -2975.6
$ws.Range("H77").Value = 4708
$ws.Range("I77").Value = 3849.6
$ws.Range("K77").Value = 19248
$ws.Range("M77").Value = -14880
$ws.Range("H110").Value = 6167901.5
$ws.Range("I110").Value = 6167901.5
$ws.Range("K110").Value = 6167901.5
$ws.Range("M110").Value = -6165856.5
$ws.Range("H122").Value = 1188.8
$ws.Range("I122").Value = 928.05884
$ws.Range("J122").Value = 2666.3333
$ws.Range("K122").Value = 2784.17652
$ws.Range("L122").Value = 7998.999899999999
$ws.Range("M122").Value = -334.17652
$ws.Range("N122").Value = -12898.9999
$ws.Range("H136").Value = 1383.8334
$ws.Range("I136").Value = 1383.8334
$ws.Range("K136").Value = 4151.5002
$ws.Range("M136").Value = -1601.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 510.7
$ws.Range("J80").Value = 444.25
$ws.Range("L80").Value = 444.25
$ws.Range("N80").Value = -2440.25
$ws.Range("H83").Value = 510.7
$ws.Range("J83").Value = 444.25
$ws.Range("L83").Value = 2221.25
$ws.Range("N83").Value = -12205.25
$ws.Range("H99").Value = 1372.6923
$ws.Range("I99").Value = 1378.75
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 1378.75
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = 119.25
$ws.Range("N99").Value = -4296
$ws.Range("H105").Value = 2708.8
$ws.Range("I105").Value = 2814.6667
$ws.Range("J105").Value = 2550
$ws.Range("K105").Value = 2814.6667
$ws.Range("L105").Value = 2550
$ws.Range("M105").Value = -1067.6667
$ws.Range("N105").Value = -6044
$ws.Range("H134").Value = 2190.7646
$ws.Range("J134").Value = 2325
$ws.Range("L134").Value = 6975
$ws.Range("N134").Value = -12045

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4300.5713
$ws.Range("I16").Value = 2787.2222
$ws.Range("K16").Value = 2787.2222
$ws.Range("M16").Value = -2500.2222
$ws.Range("H31").Value = 1744.5
$ws.Range("I31").Value = 1030.3334
$ws.Range("J31").Value = 3887
$ws.Range("K31").Value = 1030.3334
$ws.Range("L31").Value = 3887
$ws.Range("M31").Value = -735.3334
$ws.Range("N31").Value = -4477
$ws.Range("H34").Value = 1744.5
$ws.Range("I34").Value = 1030.3334
$ws.Range("J34").Value = 3887
$ws.Range("K34").Value = 1030.3334
$ws.Range("L34").Value = 3887
$ws.Range("M34").Value = -828.3334
$ws.Range("N34").Value = -4291
$ws.Range("H92").Value = 26082.666
$ws.Range("J92").Value = 26082.666
$ws.Range("L92").Value = 26082.666
$ws.Range("N92").Value = -31074.666
$ws.Range("H113").Value = 4300.5713
$ws.Range("I113").Value = 2787.2222
$ws.Range("K113").Value = 2787.2222
$ws.Range("M113").Value = -617.2222000000002
$ws.Range("H132").Value = 3397.8333
$ws.Range("I132").Value = 3467.7
$ws.Range("J132").Value = 3048.5
$ws.Range("K132").Value = 10403.1
$ws.Range("L132").Value = 9145.5
$ws.Range("M132").Value = -7873.099999999999
$ws.Range("N132").Value = -14205.5
$ws.Range("H134").Value = 2125.074
$ws.Range("I134").Value = 2259.1904
$ws.Range("J134").Value = 1655.6666
$ws.Range("K134").Value = 6777.5712
$ws.Range("L134").Value = 4966.9998
$ws.Range("M134").Value = -4242.5712
$ws.Range("N134").Value = -10036.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13750209
$ws.Range("I4").Value = 18333546
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 55000638
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -55000526
$ws.Range("N4").Value = -824
$ws.Range("H34").Value = 539
$ws.Range("I34").Value = 539
$ws.Range("K34").Value = 1617
$ws.Range("M34").Value = -1533
$ws.Range("H80").Value = 2250
$ws.Range("I80").Value = 2950
$ws.Range("J80").Value = 850
$ws.Range("K80").Value = 8850
$ws.Range("L80").Value = 2550
$ws.Range("M80").Value = -7914
$ws.Range("N80").Value = -4422
$ws.Range("H83").Value = 2250
$ws.Range("I83").Value = 2950
$ws.Range("J83").Value = 850
$ws.Range("K83").Value = 26550
$ws.Range("L83").Value = 7650
$ws.Range("M83").Value = -21870
$ws.Range("N83").Value = -17010
$ws.Range("H131").Value = 911484.4399999999
$ws.Range("J131").Value = 1002532.9
$ws.Range("L131").Value = 3007598.7
$ws.Range("N131").Value = -3017678.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 1000000000
$ws.Range("I55").Value = 1000000000
$ws.Range("K55").Value = 1000000000
$ws.Range("M55").Value = -999999673
$ws.Range("H122").Value = 784.8333
$ws.Range("I122").Value = 784.8333
$ws.Range("K122").Value = 2354.4999
$ws.Range("M122").Value = 95.5001000000002
$ws.Range("H126").Value = 2249.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 906.8570999999999
$ws.Range("I22").Value = 892.1667
$ws.Range("J22").Value = 995
$ws.Range("K22").Value = 892.1667
$ws.Range("L22").Value = 995
$ws.Range("M22").Value = -597.1667
$ws.Range("N22").Value = -1585
$ws.Range("H27").Value = 906.8570999999999
$ws.Range("I27").Value = 892.1667
$ws.Range("J27").Value = 995
$ws.Range("K27").Value = 892.1667
$ws.Range("L27").Value = 995
$ws.Range("M27").Value = -785.1667
$ws.Range("N27").Value = -1209
$ws.Range("H55").Value = 1042.8889
$ws.Range("I55").Value = 414.66666
$ws.Range("J55").Value = 2299.3333
$ws.Range("K55").Value = 414.66666
$ws.Range("L55").Value = 2299.3333
$ws.Range("M55").Value = -241.66666
$ws.Range("N55").Value = -2645.3333
$ws.Range("H122").Value = 8424.639999999999
$ws.Range("I122").Value = 8187.5
$ws.Range("K122").Value = 24562.5
$ws.Range("M122").Value = -22112.5
$ws.Range("H136").Value = 4799
$ws.Range("I136").Value = 4332.6665
$ws.Range("J136").Value = 5498.5
$ws.Range("K136").Value = 12997.9995
$ws.Range("L136").Value = 16495.5
$ws.Range("M136").Value = -10447.9995
$ws.Range("N136").Value = -21595.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 13772499
$ws.Range("I4").Value = 25002500
$ws.Range("J4").Value = 2542497.5
$ws.Range("K4").Value = 25002500
$ws.Range("L4").Value = 2542497.5
$ws.Range("M4").Value = -25002387
$ws.Range("N4").Value = -2542723.5
$ws.Range("H56").Value = 26799.334
$ws.Range("J56").Value = 26799.334
$ws.Range("L56").Value = 26799.334
$ws.Range("N56").Value = -28227.334
$ws.Range("H122").Value = 1754.8889
$ws.Range("I122").Value = 1466.6666
$ws.Range("J122").Value = 2331.3333
$ws.Range("K122").Value = 4399.9998
$ws.Range("L122").Value = 6993.999899999999
$ws.Range("M122").Value = -1949.9998
$ws.Range("N122").Value = -11893.9999
$ws.Range("H126").Value = 4366
$ws.Range("I126").Value = 4308.6
$ws.Range("K126").Value = 12925.8
$ws.Range("M126").Value = -10455.8
